$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 data
$ws.Range("A10").Value = 42352
$ws.Range("B10").Value = "cbardash"

# NOTE: shared-string creation order matters (new strings are appended to
# sharedStrings.xml in the order they are first referenced), so write the
# long "Description of Updates" (D) before the short "Commit Comment" (C)
# to land them at the same indices as the target workbook.
$ws.Range("D10").Value = "Added survey popup screen on startup, and also survey button at bottom of TOC.`nRemoved esriAttribution from lower right-hand corner.`nSwitched ArcGIS Server-based services over to AGO-hosted services.`nAdded version number to About tab."
$ws.Range("C10").Value = "added survey/removed attribution from lower right-hand corner/removed ArcGIS Server layer URLs/added version number"

$ws.Range("E10").Value = "see SPM_TestScrip.docx"
$ws.Range("F10").Value = "YES"

$ws.Rows.Item(10).RowHeight = 60

$ws.Range("A7").Select()
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("C10").Select()
